$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2172.6086
$ws.Range("J17").Value = 2172.6086
$ws.Range("L17").Value = 6517.825800000001
$ws.Range("N17").Value = -6853.825800000001

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9492.643
$ws.Range("I51").Value = 7400
$ws.Range("J51").Value = 9841.416999999999
$ws.Range("K51").Value = 7400
$ws.Range("L51").Value = 9841.416999999999
$ws.Range("M51").Value = -6916
$ws.Range("N51").Value = -10809.417

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5164.125
$ws.Range("I62").Value = 4973.2856
$ws.Range("K62").Value = 4973.2856
$ws.Range("M62").Value = -4349.2856

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6599.9
$ws.Range("J64").Value = 5125
$ws.Range("L64").Value = 5125
$ws.Range("N64").Value = -5621

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5164.125
$ws.Range("I65").Value = 4973.2856
$ws.Range("K65").Value = 24866.428
$ws.Range("M65").Value = -21746.428

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6599.9
$ws.Range("J67").Value = 5125
$ws.Range("L67").Value = 5125
$ws.Range("N67").Value = -6841

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6338.8667
$ws.Range("I74").Value = 5929.4614
$ws.Range("K74").Value = 5929.4614
$ws.Range("M74").Value = -4993.4614

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 6338.8667
$ws.Range("I77").Value = 5929.4614
$ws.Range("K77").Value = 29647.307
$ws.Range("M77").Value = -24967.307

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I80").Value = 1149.75
$ws.Range("J80").Value = 678.63635
$ws.Range("K80").Value = 3449.25
$ws.Range("L80").Value = 2035.90905
$ws.Range("M80").Value = -2451.25
$ws.Range("N80").Value = -4031.90905

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I83").Value = 1149.75
$ws.Range("J83").Value = 678.63635
$ws.Range("K83").Value = 10347.75
$ws.Range("L83").Value = 6107.72715
$ws.Range("M83").Value = -5355.75
$ws.Range("N83").Value = -16091.72715

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4206.727
$ws.Range("I88").Value = 5399.6
$ws.Range("J88").Value = 3212.6667
$ws.Range("K88").Value = 5399.6
$ws.Range("L88").Value = 3212.6667
$ws.Range("M88").Value = -4993.6
$ws.Range("N88").Value = -4024.6667

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4206.727
$ws.Range("I91").Value = 5399.6
$ws.Range("J91").Value = 3212.6667
$ws.Range("K91").Value = 5399.6
$ws.Range("L91").Value = 3212.6667
$ws.Range("M91").Value = -3995.6
$ws.Range("N91").Value = -6020.6667

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4459.85
$ws.Range("I113").Value = 3531.3845
$ws.Range("J113").Value = 6184.143
$ws.Range("K113").Value = 3531.3845
$ws.Range("L113").Value = 6184.143
$ws.Range("M113").Value = -277.3845000000001
$ws.Range("N113").Value = -12692.143

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11110.263
$ws.Range("I116").Value = 9522.333000000001
$ws.Range("J116").Value = 12539.4
$ws.Range("K116").Value = 9522.333000000001
$ws.Range("L116").Value = 12539.4
$ws.Range("M116").Value = -6080.333000000001
$ws.Range("N116").Value = -19423.4

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 5002
$ws.Range("J121").Value = 5002
$ws.Range("L121").Value = 15006
$ws.Range("N121").Value = -18500

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1499.5
$ws.Range("I131").Value = 1499.5
$ws.Range("K131").Value = 4498.5
$ws.Range("M131").Value = 541.5

# ALC row 134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 82064.836
$ws.Range("J134").Value = 82064.836
$ws.Range("L134").Value = 82064.836
$ws.Range("N134").Value = -92204.836

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1607.1578
$ws.Range("I137").Value = 1359.3143
$ws.Range("K137").Value = 4077.9429
$ws.Range("M137").Value = -1527.9429

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2927.5
$ws.Range("J138").Value = 3306.0488
$ws.Range("L138").Value = 9918.1464
$ws.Range("N138").Value = -20198.1464

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 94497.8
$ws.Range("J140").Value = 94497.8
$ws.Range("L140").Value = 94497.8
$ws.Range("N140").Value = -104857.8

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 191.11111
$ws.Range("I5").Value = 231.5
$ws.Range("K5").Value = 231.5
$ws.Range("M5").Value = -119.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4776.1626
$ws.Range("I32").Value = 2235.948
$ws.Range("J32").Value = 26509.111
$ws.Range("K32").Value = 2235.948
$ws.Range("L32").Value = 26509.111
$ws.Range("M32").Value = -1948.948
$ws.Range("N32").Value = -27083.111

# ARM row 59
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10311.111
$ws.Range("I61").Value = 10975
$ws.Range("K61").Value = 10975
$ws.Range("M61").Value = -10763

# ARM row 64
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 50000
$ws.Range("K64").Value = 50000
$ws.Range("M64").Value = -49752

# ARM row 67
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 50000
$ws.Range("K67").Value = 50000
$ws.Range("M67").Value = -49142

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4250
$ws.Range("I74").Value = 1735.85
$ws.Range("J74").Value = 8117.923
$ws.Range("K74").Value = 1735.85
$ws.Range("L74").Value = 8117.923
$ws.Range("M74").Value = -861.8499999999999
$ws.Range("N74").Value = -9865.922999999999

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4250
$ws.Range("I77").Value = 1735.85
$ws.Range("J77").Value = 8117.923
$ws.Range("K77").Value = 8679.25
$ws.Range("L77").Value = 40589.615
$ws.Range("M77").Value = -4311.25
$ws.Range("N77").Value = -49325.615

# ARM row 109
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 47498.25
$ws.Range("J109").Value = 47498.25
$ws.Range("L109").Value = 47498.25
$ws.Range("N109").Value = -50272.25

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3173.9375
$ws.Range("I122").Value = 2076.5557
$ws.Range("K122").Value = 6229.6671
$ws.Range("M122").Value = -3779.6671

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5435.0415
$ws.Range("I132").Value = 2680.6155
$ws.Range("K132").Value = 8041.8465
$ws.Range("M132").Value = -5511.8465

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10311.111
$ws.Range("I136").Value = 10975
$ws.Range("K136").Value = 32925
$ws.Range("M136").Value = -30375

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 191.11111
$ws.Range("I4").Value = 231.5
$ws.Range("K4").Value = 231.5
$ws.Range("M4").Value = -116.5

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2295.3845
$ws.Range("I20").Value = 1459.1111
$ws.Range("K20").Value = 1459.1111
$ws.Range("M20").Value = -1212.1111

# BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 9990
$ws.Range("I62").Value = 9990
$ws.Range("K62").Value = 9990
$ws.Range("M62").Value = -9304

# BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 9990
$ws.Range("I65").Value = 9990
$ws.Range("K65").Value = 29970
$ws.Range("M65").Value = -26538

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 959.34375
$ws.Range("I94").Value = 708.9655
$ws.Range("K94").Value = 708.9655
$ws.Range("M94").Value = -257.9655

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2869.8235
$ws.Range("I99").Value = 2767.9375
$ws.Range("K99").Value = 2767.9375
$ws.Range("M99").Value = -1269.9375

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3746.8333
$ws.Range("I105").Value = 3840.25
$ws.Range("J105").Value = 3560
$ws.Range("K105").Value = 3840.25
$ws.Range("L105").Value = 3560
$ws.Range("M105").Value = -2093.25
$ws.Range("N105").Value = -7054

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3312.4
$ws.Range("I107").Value = 3773.5
$ws.Range("J107").Value = 2725.5454
$ws.Range("K107").Value = 3773.5
$ws.Range("L107").Value = 2725.5454
$ws.Range("M107").Value = -1853.5
$ws.Range("N107").Value = -6565.5454

# BSM row 108
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 99000
$ws.Range("J108").Value = 99000
$ws.Range("L108").Value = 99000
$ws.Range("N108").Value = -106680

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3822.5417
$ws.Range("I134").Value = 2394.55
$ws.Range("J134").Value = 10962.5
$ws.Range("K134").Value = 7183.650000000001
$ws.Range("L134").Value = 32887.5
$ws.Range("M134").Value = -4648.650000000001
$ws.Range("N134").Value = -37957.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6266.6
$ws.Range("I31").Value = 2906.3809
$ws.Range("K31").Value = 2906.3809
$ws.Range("M31").Value = -2611.3809

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6266.6
$ws.Range("I34").Value = 2906.3809
$ws.Range("K34").Value = 2906.3809
$ws.Range("M34").Value = -2704.3809

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3362.5908
$ws.Range("I58").Value = 1998.9333
$ws.Range("J58").Value = 6284.7144
$ws.Range("K58").Value = 1998.9333
$ws.Range("L58").Value = 6284.7144
$ws.Range("M58").Value = -1795.9333
$ws.Range("N58").Value = -6690.7144

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3367.611
$ws.Range("I62").Value = 2998
$ws.Range("J62").Value = 3948.4285
$ws.Range("K62").Value = 2998
$ws.Range("L62").Value = 3948.4285
$ws.Range("M62").Value = -2374
$ws.Range("N62").Value = -5196.4285

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3367.611
$ws.Range("I65").Value = 2998
$ws.Range("J65").Value = 3948.4285
$ws.Range("K65").Value = 14990
$ws.Range("L65").Value = 19742.1425
$ws.Range("M65").Value = -11870
$ws.Range("N65").Value = -25982.1425

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5791.7856
$ws.Range("I86").Value = 6023.75
$ws.Range("J86").Value = 5482.5
$ws.Range("K86").Value = 6023.75
$ws.Range("L86").Value = 5482.5
$ws.Range("M86").Value = -4900.75
$ws.Range("N86").Value = -7728.5

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5791.7856
$ws.Range("I89").Value = 6023.75
$ws.Range("J89").Value = 5482.5
$ws.Range("K89").Value = 30118.75
$ws.Range("L89").Value = 27412.5
$ws.Range("M89").Value = -24502.75
$ws.Range("N89").Value = -38644.5

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2508.3125
$ws.Range("I122").Value = 2291.6667
$ws.Range("K122").Value = 6875.000100000001
$ws.Range("M122").Value = -4425.000100000001

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5101.125
$ws.Range("I132").Value = 2801.75
$ws.Range("K132").Value = 8405.25
$ws.Range("M132").Value = -5875.25

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5000.385
$ws.Range("I134").Value = 4221.3335
$ws.Range("J134").Value = 6753.25
$ws.Range("K134").Value = 12664.0005
$ws.Range("L134").Value = 20259.75
$ws.Range("M134").Value = -10129.0005
$ws.Range("N134").Value = -25329.75

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3362.5908
$ws.Range("I136").Value = 1998.9333
$ws.Range("J136").Value = 6284.7144
$ws.Range("K136").Value = 5996.7999
$ws.Range("L136").Value = 18854.1432
$ws.Range("M136").Value = -3446.7999
$ws.Range("N136").Value = -23954.1432

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 7.5
$ws.Range("I7").Value = 7.5
$ws.Range("K7").Value = 22.5
$ws.Range("M7").Value = 89.5

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3924
$ws.Range("I80").Value = 4124.75
$ws.Range("J80").Value = 3723.25
$ws.Range("K80").Value = 12374.25
$ws.Range("L80").Value = 11169.75
$ws.Range("M80").Value = -11438.25
$ws.Range("N80").Value = -13041.75

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3924
$ws.Range("I83").Value = 4124.75
$ws.Range("J83").Value = 3723.25
$ws.Range("K83").Value = 37122.75
$ws.Range("L83").Value = 33509.25
$ws.Range("M83").Value = -32442.75
$ws.Range("N83").Value = -42869.25

# CUL row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 12398.4
$ws.Range("I120").Value = 9595.799999999999
$ws.Range("J120").Value = 13799.7
$ws.Range("K120").Value = 28787.4
$ws.Range("L120").Value = 41399.10000000001
$ws.Range("M120").Value = -23949.4
$ws.Range("N120").Value = -51075.10000000001

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 959976.5
$ws.Range("J131").Value = 1823020.6
$ws.Range("L131").Value = 5469061.800000001
$ws.Range("N131").Value = -5479141.800000001

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3626.0476
$ws.Range("I141").Value = 2848
$ws.Range("K141").Value = 8544
$ws.Range("M141").Value = -3364

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3757.1304
$ws.Range("I102").Value = 2946.6843
$ws.Range("J102").Value = 7606.75
$ws.Range("K102").Value = 2946.6843
$ws.Range("L102").Value = 7606.75
$ws.Range("M102").Value = -1324.6843
$ws.Range("N102").Value = -10850.75

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9585.214
$ws.Range("I122").Value = 6170.7144
$ws.Range("J122").Value = 12999.714
$ws.Range("K122").Value = 18512.1432
$ws.Range("L122").Value = 38999.142
$ws.Range("M122").Value = -16062.1432
$ws.Range("N122").Value = -43899.142

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4736.6855
$ws.Range("I132").Value = 3070.8
$ws.Range("J132").Value = 6957.8667
$ws.Range("K132").Value = 9212.400000000001
$ws.Range("L132").Value = 20873.6001
$ws.Range("M132").Value = -6682.400000000001
$ws.Range("N132").Value = -25933.6001

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7288.579
$ws.Range("I122").Value = 6217.6875
$ws.Range("J122").Value = 13000
$ws.Range("K122").Value = 18653.0625
$ws.Range("L122").Value = 39000
$ws.Range("M122").Value = -16203.0625
$ws.Range("N122").Value = -43900

# LTW row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 66476
$ws.Range("J128").Value = 66476
$ws.Range("L128").Value = 66476
$ws.Range("N128").Value = -76436

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5537.154
$ws.Range("I132").Value = 2961.6667
$ws.Range("K132").Value = 8885.000100000001
$ws.Range("M132").Value = -6355.000100000001

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3385.3125
$ws.Range("I81").Value = 3337.6667
$ws.Range("K81").Value = 6675.3334
$ws.Range("M81").Value = -5614.3334

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3385.3125
$ws.Range("I84").Value = 3337.6667
$ws.Range("K84").Value = 33376.667
$ws.Range("M84").Value = -28072.667

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 757.4
$ws.Range("I100").Value = 499.66666
$ws.Range("K100").Value = 999.33332
$ws.Range("M100").Value = -458.33332

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3005.138
$ws.Range("I122").Value = 2920.6316
$ws.Range("J122").Value = 3165.7
$ws.Range("K122").Value = 8761.8948
$ws.Range("L122").Value = 9497.099999999999
$ws.Range("M122").Value = -6311.8948
$ws.Range("N122").Value = -14397.1

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15685.167
$ws.Range("I136").Value = 41111
$ws.Range("K136").Value = 123333
$ws.Range("M136").Value = -120783
